$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text (locale-formatted, e.g. thousand-dot separators)
# that must stay literal text rather than being auto-parsed into a float by
# Excel, so force a Text number format on each D cell before writing its value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.295.20'
$ws.Range("E2").Value = '  -0.47%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.840.08'
$ws.Range("E3").Value = '  -0.53%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9971'
$ws.Range("E4").Value = '  -0.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.23'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6271'
$ws.Range("E6").Value = '  -0.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9993'
$ws.Range("E7").Value = '  -0.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07486'
$ws.Range("E8").Value = '  -2.55%  '

$ws.Range("E9").Value = '  -0.91%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07711'
$ws.Range("E11").Value = '  -0.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.841.42'
$ws.Range("E12").Value = '  -0.70%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.987'
$ws.Range("E13").Value = '  -1.04%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6779'
$ws.Range("E14").Value = '  -0.56%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001028'
$ws.Range("E15").Value = '  -4.51%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.05'
$ws.Range("E16").Value = '  -1.77%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.100.74'
$ws.Range("E17").Value = '  -0.39%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.104'
$ws.Range("E18").Value = '  -1.84%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.313.11'
$ws.Range("E19").Value = '  -0.49%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '228.57'
$ws.Range("E20").Value = '  -0.11%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.26'
$ws.Range("E21").Value = '  -1.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9989'
$ws.Range("E22").Value = '  -0.18%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.379'
$ws.Range("E23").Value = '  -1.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9993'
$ws.Range("E24").Value = '  -0.17%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.61'
$ws.Range("E25").Value = '  +0.59%  '

$ws.Range("E26").Value = '  +0.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.381'
$ws.Range("E27").Value = '  -0.43%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.52'
$ws.Range("E28").Value = '  -1.22%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.393'
$ws.Range("E29").Value = '  +1.87%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.471'
$ws.Range("E30").Value = '  +0.73%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05703'
$ws.Range("E31").Value = '  +1.24%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.095'
$ws.Range("E32").Value = '  -0.80%  '

$ws.Range("E33").Value = '  -0.68%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.819'
$ws.Range("E34").Value = '  -1.46%  '

$ws.Range("E35").Value = '  -1.82%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6925'
$ws.Range("E36").Value = '  -2.04%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.583'
$ws.Range("E37").Value = '  -0.51%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.818'
$ws.Range("E38").Value = '  +2.26%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.242.51'
$ws.Range("E39").Value = '  +1.34%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01814'
$ws.Range("E40").Value = '  +1.09%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.510'
$ws.Range("E41").Value = '  +0.72%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9052'
$ws.Range("E42").Value = '  +0.19%  '

$ws.Range("E43").Value = '  -0.22%  '

$ws.Range("E44").Value = '  -0.63%  '

$ws.Range("E45").Value = '  -0.76%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.75'
$ws.Range("E46").Value = '  -0.61%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.072'
$ws.Range("E47").Value = '  -1.60%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1164'
$ws.Range("E48").Value = '  +0.63%  '

$ws.Range("E49").Value = '  -2.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.961'
$ws.Range("E50").Value = '  -0.98%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3934'
$ws.Range("E51").Value = '  -2.20%  '
